$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.654.85'
$ws.Range("E2").Value = '  +3.50%  '
$ws.Range("D3").Value = '1.610.00'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = "'212.76"
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("E6").Value = '  +1.89%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = "'27.14"
$ws.Range("E8").Value = '  +9.22%  '
$ws.Range("D9").Value = "'43.65"
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  +2.44%  '
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("D12").Value = "'0.0911"
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").Value = '1.840.61'
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("D14").Value = '1.604.52'
$ws.Range("E14").Value = '  +2.52%  '
$ws.Range("D15").Value = '29.665.95'
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("D18").Value = "'63.57"
$ws.Range("E18").Value = '  +3.41%  '
$ws.Range("D19").Value = "'240.51"
$ws.Range("E19").Value = '  +5.56%  '
$ws.Range("E20").Value = '  +3.79%  '
$ws.Range("D21").Value = '0.0₃0695'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = "'3.99"
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("D24").Value = "'9.24"
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("D26").Value = "'154.90"
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("D27").Value = "'15.34"
$ws.Range("E27").Value = '  +3.85%  '
$ws.Range("E28").Value = '  +2.23%  '
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  +3.78%  '
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").Value = "'3.22"
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").Value = "'3.14"
$ws.Range("E34").Value = '  +4.51%  '
$ws.Range("D35").Value = '1.428.06'
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  +4.77%  '
$ws.Range("E38").Value = '  +5.71%  '
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("E41").Value = '  +4.52%  '
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("D43").Value = "'54.56"
$ws.Range("E43").Value = '  +27.56%  '
$ws.Range("E44").Value = '  +6.20%  '
$ws.Range("E45").Value = '  +4.35%  '
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = "'65.97"
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").Value = "'5.30"
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").Value = '1.750.89'
$ws.Range("E49").Value = '  +3.07%  '
$ws.Range("D50").Value = "'0.902"
$ws.Range("E50").Value = '  +5.19%  '
$ws.Range("D51").Value = "'86.83"
$ws.Range("E51").Value = '  +2.26%  '
